$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sector_ID"
$ws.Range("B2").Value = "BusinessKey"
$ws.Range("C2").Value = "ProgrammeBusinessKey"
$ws.Range("D2").Value = "Code"
$ws.Range("E2").Value = "LongName"
$ws.Range("F2").Value = "ShortName"
$ws.Range("G2").Value = "TextDescription"
